# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gains a new (blank) column between the
# existing "Late" column (N) and the "Outstanding" / heading / "Outstanding"
# columns that follow it (O, P). Inserting a whole column there shifts the
# old N/O/P data right to O/P/Q and leaves a new, empty N column behind -
# exactly what the canonical XML shows.
#
# It also records that the "Repayment schedule" tab (not "Edit Repayment
# Schedule") is the one the user last had selected/active.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Match the new column's width to its neighbour (column M / "In Advance")
# before inserting, then re-apply the same width to the freshly inserted
# column N.
$existingWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column at N; this shifts old N->O, O->P, P->Q.
$ws.Columns("N").Insert()

$ws.Columns("N").ColumnWidth = $existingWidth

# Make "Repayment schedule" the active/selected sheet and cell, matching
# the saved view state in the workbook.
$ws.Activate()
[void]$ws.Range("T9").Select()
